# Apply updated loading_percent values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 8.548047572575818
$ws.Cells.Item(2, 4).Value = 6.54125883045493
$ws.Cells.Item(2, 5).Value = 12.18932593745597
$ws.Cells.Item(2, 6).Value = 38.76473174275738
$ws.Cells.Item(2, 7).Value = 56.05567591988871
$ws.Cells.Item(2, 8).Value = 19.44283405269406
$ws.Cells.Item(2, 9).Value = 33.1090684750811
$ws.Cells.Item(2, 12).Value = 9.605220242882643
$ws.Cells.Item(2, 14).Value = 18.99769123883691

# Row 3
$ws.Cells.Item(3, 3).Value = 8.564959229818934
$ws.Cells.Item(3, 4).Value = 6.524013138894491
$ws.Cells.Item(3, 5).Value = 12.17799175561252
$ws.Cells.Item(3, 6).Value = 37.95564554768515
$ws.Cells.Item(3, 7).Value = 54.44485957427088
$ws.Cells.Item(3, 8).Value = 19.25349562898403
$ws.Cells.Item(3, 9).Value = 32.61190476574679
$ws.Cells.Item(3, 12).Value = 9.609539774742887
$ws.Cells.Item(3, 14).Value = 18.40031508502705

# Row 4
$ws.Cells.Item(4, 3).Value = 8.57669863710624
$ws.Cells.Item(4, 4).Value = 6.51520600057889
$ws.Cells.Item(4, 5).Value = 12.17398774936703
$ws.Cells.Item(4, 6).Value = 37.46450684535414
$ws.Cells.Item(4, 7).Value = 53.45322415713718
$ws.Cells.Item(4, 8).Value = 19.14253674959212
$ws.Cells.Item(4, 9).Value = 32.31443100047505
$ws.Cells.Item(4, 12).Value = 9.614265052071493
$ws.Cells.Item(4, 14).Value = 18.02485520896364

# Row 5
$ws.Cells.Item(5, 3).Value = 8.581822926415029
$ws.Cells.Item(5, 4).Value = 6.51206511652173
$ws.Cells.Item(5, 5).Value = 12.17309801068173
$ws.Cells.Item(5, 6).Value = 37.26607771508598
$ws.Cells.Item(5, 7).Value = 53.04915829273452
$ws.Cells.Item(5, 8).Value = 19.09869178991451
$ws.Cells.Item(5, 9).Value = 32.19531226839234
$ws.Cells.Item(5, 12).Value = 9.616710417996153
$ws.Cells.Item(5, 14).Value = 17.86990355188116

# Row 6
$ws.Cells.Item(6, 3).Value = 8.582694354639402
$ws.Cells.Item(6, 4).Value = 6.511570631173162
$ws.Cells.Item(6, 5).Value = 12.17299503610342
$ws.Cells.Item(6, 6).Value = 37.23324084759345
$ws.Cells.Item(6, 7).Value = 52.98208606689533
$ws.Cells.Item(6, 8).Value = 19.09149532236653
$ws.Cells.Item(6, 9).Value = 32.17566395987912
$ws.Cells.Item(6, 12).Value = 9.617147814398114
$ws.Cells.Item(6, 14).Value = 17.84406337566596

# Row 7
$ws.Cells.Item(7, 3).Value = 8.576766367500783
$ws.Cells.Item(7, 4).Value = 6.515161827466969
$ws.Cells.Item(7, 5).Value = 12.17397274793382
$ws.Cells.Item(7, 6).Value = 37.46182343124732
$ws.Cells.Item(7, 7).Value = 53.44777372132972
$ws.Cells.Item(7, 8).Value = 19.14193983662412
$ws.Cells.Item(7, 9).Value = 32.31281581399847
$ws.Cells.Item(7, 12).Value = 9.614295928700436
$ws.Cells.Item(7, 14).Value = 18.02277304766463

# Row 8
$ws.Cells.Item(8, 3).Value = 8.553597072834577
$ws.Cells.Item(8, 4).Value = 6.53494287687809
$ws.Cells.Item(8, 5).Value = 12.18480398446885
$ws.Cells.Item(8, 6).Value = 38.48476391244221
$ws.Cells.Item(8, 7).Value = 55.50120440250834
$ws.Cells.Item(8, 8).Value = 19.37647600574521
$ws.Cells.Item(8, 9).Value = 32.93611938886439
$ws.Cells.Item(8, 12).Value = 9.60627849710111
$ws.Cells.Item(8, 14).Value = 18.79364780656866

# Row 9
$ws.Cells.Item(9, 3).Value = 8.51893989266928
$ws.Cells.Item(9, 4).Value = 6.587848814916436
$ws.Cells.Item(9, 5).Value = 12.22952545578663
$ws.Cells.Item(9, 6).Value = 40.52271384069176
$ws.Cells.Item(9, 7).Value = 59.47908260380215
$ws.Cells.Item(9, 8).Value = 19.87651855972801
$ws.Cells.Item(9, 9).Value = 34.21358298927608
$ws.Cells.Item(9, 12).Value = 9.607065939956465
$ws.Cells.Item(9, 14).Value = 20.2273683202997

# Row 10
$ws.Cells.Item(10, 3).Value = 8.500077603047993
$ws.Cells.Item(10, 4).Value = 6.635265484924069
$ws.Cells.Item(10, 5).Value = 12.27671356192986
$ws.Cells.Item(10, 6).Value = 42.02276729776702
$ws.Cells.Item(10, 7).Value = 62.33661497587624
$ws.Cells.Item(10, 8).Value = 20.26563986766339
$ws.Cells.Item(10, 9).Value = 35.17675798586905
$ws.Cells.Item(10, 12).Value = 9.617782190532838
$ws.Cells.Item(10, 14).Value = 21.22223697909767

# Row 11
$ws.Cells.Item(11, 3).Value = 8.492935977384757
$ws.Cells.Item(11, 4).Value = 6.658670252532967
$ws.Cells.Item(11, 5).Value = 12.30128531533712
$ws.Cells.Item(11, 6).Value = 42.70243386384655
$ws.Cells.Item(11, 7).Value = 63.6158804261721
$ws.Cells.Item(11, 8).Value = 20.44676437197159
$ws.Cells.Item(11, 9).Value = 35.6183635271024
$ws.Cells.Item(11, 12).Value = 9.624870026512118
$ws.Cells.Item(11, 14).Value = 21.66018057919901

# Row 12
$ws.Cells.Item(12, 3).Value = 8.490439001446719
$ws.Cells.Item(12, 4).Value = 6.667794205637503
$ws.Cells.Item(12, 5).Value = 12.31103500637909
$ws.Cells.Item(12, 6).Value = 42.95915356060193
$ws.Cells.Item(12, 7).Value = 64.09684983928501
$ws.Cells.Item(12, 8).Value = 20.51588900938445
$ws.Cells.Item(12, 9).Value = 35.78592729203418
$ws.Cells.Item(12, 12).Value = 9.627872742732462
$ws.Cells.Item(12, 14).Value = 21.82377585682186

# Row 13
$ws.Cells.Item(13, 3).Value = 8.490967538079476
$ws.Cells.Item(13, 4).Value = 6.665817636079578
$ws.Cells.Item(13, 5).Value = 12.30891548132112
$ws.Cells.Item(13, 6).Value = 42.90389758222634
$ws.Cells.Item(13, 7).Value = 63.99342574258328
$ws.Cells.Item(13, 8).Value = 20.50097878508261
$ws.Cells.Item(13, 9).Value = 35.74982692619119
$ws.Cells.Item(13, 12).Value = 9.627211874153534
$ws.Cells.Item(13, 14).Value = 21.78864458690801

# Row 14
$ws.Cells.Item(14, 3).Value = 8.492726390302192
$ws.Cells.Item(14, 4).Value = 6.659415670601665
$ws.Cells.Item(14, 5).Value = 12.30207852333074
$ws.Cells.Item(14, 6).Value = 42.72356894458952
$ws.Cells.Item(14, 7).Value = 63.65552208521829
$ws.Cells.Item(14, 8).Value = 20.45244085516624
$ws.Cells.Item(14, 9).Value = 35.63214323595626
$ws.Cells.Item(14, 12).Value = 9.625110672440348
$ws.Cells.Item(14, 14).Value = 21.67368539489659

# Row 15
$ws.Cells.Item(15, 3).Value = 8.493830762105409
$ws.Cells.Item(15, 4).Value = 6.655528192578607
$ws.Cells.Item(15, 5).Value = 12.29794856939128
$ws.Cells.Item(15, 6).Value = 42.61301950385379
$ws.Cells.Item(15, 7).Value = 63.44808181910363
$ws.Cells.Item(15, 8).Value = 20.42277824019309
$ws.Cells.Item(15, 9).Value = 35.56009785080179
$ws.Cells.Item(15, 12).Value = 9.623865141324737
$ws.Cells.Item(15, 14).Value = 21.60297336126124

# Row 16
$ws.Cells.Item(16, 3).Value = 8.500573322065993
$ws.Cells.Item(16, 4).Value = 6.633772621830958
$ws.Cells.Item(16, 5).Value = 12.27517010146205
$ws.Cells.Item(16, 6).Value = 41.97827187745378
$ws.Cells.Item(16, 7).Value = 62.25255348183828
$ws.Cells.Item(16, 8).Value = 20.25388104402256
$ws.Cells.Item(16, 9).Value = 35.14795381053022
$ws.Cells.Item(16, 12).Value = 9.617363560804215
$ws.Cells.Item(16, 14).Value = 21.19330956972086

# Row 17
$ws.Cells.Item(17, 3).Value = 8.50507856087312
$ws.Cells.Item(17, 4).Value = 6.620894374435061
$ws.Cells.Item(17, 5).Value = 12.2619906475488
$ws.Cells.Item(17, 6).Value = 41.58798627937279
$ws.Cells.Item(17, 7).Value = 61.51348407508949
$ws.Cells.Item(17, 8).Value = 20.15128256180353
$ws.Cells.Item(17, 9).Value = 34.89588669837362
$ws.Cells.Item(17, 12).Value = 9.613942260222533
$ws.Cells.Item(17, 14).Value = 20.93814219015166

# Row 18
$ws.Cells.Item(18, 3).Value = 8.507805265770457
$ws.Cells.Item(18, 4).Value = 6.613660029825692
$ws.Cells.Item(18, 5).Value = 12.25470262740408
$ws.Cells.Item(18, 6).Value = 41.36326370541828
$ws.Cells.Item(18, 7).Value = 61.08647466578538
$ws.Cells.Item(18, 8).Value = 20.09266092196951
$ws.Cells.Item(18, 9).Value = 34.75123864885476
$ws.Cells.Item(18, 12).Value = 9.612182698080019
$ws.Cells.Item(18, 14).Value = 20.79000725568362

# Row 19
$ws.Cells.Item(19, 3).Value = 8.508751722894061
$ws.Cells.Item(19, 4).Value = 6.6112403743689
$ws.Cells.Item(19, 5).Value = 12.2522853008525
$ws.Cells.Item(19, 6).Value = 41.28714348487451
$ws.Cells.Item(19, 7).Value = 60.94158326043432
$ws.Cells.Item(19, 8).Value = 20.0728813228323
$ws.Cells.Item(19, 9).Value = 34.70232562587087
$ws.Cells.Item(19, 12).Value = 9.611622694298317
$ws.Cells.Item(19, 14).Value = 20.73962067985785

# Row 20
$ws.Cells.Item(20, 3).Value = 8.504584952423636
$ws.Cells.Item(20, 4).Value = 6.622247411689989
$ws.Cells.Item(20, 5).Value = 12.26336336567635
$ws.Cells.Item(20, 6).Value = 41.62955971611812
$ws.Cells.Item(20, 7).Value = 61.59236131146599
$ws.Cells.Item(20, 8).Value = 20.1621643439278
$ws.Cells.Item(20, 9).Value = 34.9226862398891
$ws.Cells.Item(20, 12).Value = 9.614284903154781
$ws.Cells.Item(20, 14).Value = 20.96544799483449

# Row 21
$ws.Cells.Item(21, 3).Value = 8.492204140146491
$ws.Cells.Item(21, 4).Value = 6.661289021603446
$ws.Cells.Item(21, 5).Value = 12.30407464518058
$ws.Cells.Item(21, 6).Value = 42.77655565751448
$ws.Cells.Item(21, 7).Value = 63.75487020094825
$ws.Cells.Item(21, 8).Value = 20.46668349264012
$ws.Cells.Item(21, 9).Value = 35.66670188230873
$ws.Cells.Item(21, 12).Value = 9.625719193949806
$ws.Cells.Item(21, 14).Value = 21.70751365554064

# Row 22
$ws.Cells.Item(22, 3).Value = 8.485321675577307
$ws.Cells.Item(22, 4).Value = 6.688325141210989
$ws.Cells.Item(22, 5).Value = 12.33327389423361
$ws.Cells.Item(22, 6).Value = 43.52225221846174
$ws.Cells.Item(22, 7).Value = 65.14782695422065
$ws.Cells.Item(22, 8).Value = 20.66880714239005
$ws.Cells.Item(22, 9).Value = 36.15485449008164
$ws.Cells.Item(22, 12).Value = 9.635049897478194
$ws.Cells.Item(22, 14).Value = 22.17935961385674

# Row 23
$ws.Cells.Item(23, 3).Value = 8.488884187410376
$ws.Cells.Item(23, 4).Value = 6.673757379353954
$ws.Cells.Item(23, 5).Value = 12.31745322297694
$ws.Cells.Item(23, 6).Value = 43.12470326521412
$ws.Cells.Item(23, 7).Value = 64.40639437311168
$ws.Cells.Item(23, 8).Value = 20.56066418191492
$ws.Cells.Item(23, 9).Value = 35.89419634587594
$ws.Cells.Item(23, 12).Value = 9.629899839867083
$ws.Cells.Item(23, 14).Value = 21.92877110911181

# Row 24
$ws.Cells.Item(24, 3).Value = 8.504807687414248
$ws.Cells.Item(24, 4).Value = 6.621635175291531
$ws.Cells.Item(24, 5).Value = 12.26274185944519
$ws.Cells.Item(24, 6).Value = 41.61076541093709
$ws.Cells.Item(24, 7).Value = 61.55670744272783
$ws.Cells.Item(24, 8).Value = 20.15724355359084
$ws.Cells.Item(24, 9).Value = 34.91056931814398
$ws.Cells.Item(24, 12).Value = 9.614129348353503
$ws.Cells.Item(24, 14).Value = 20.95310750188673

# Row 25
$ws.Cells.Item(25, 3).Value = 8.527158575341213
$ws.Cells.Item(25, 4).Value = 6.572029351468887
$ws.Cells.Item(25, 5).Value = 12.21490947342395
$ws.Cells.Item(25, 6).Value = 39.96973144488514
$ws.Cells.Item(25, 7).Value = 58.41199414297472
$ws.Cells.Item(25, 8).Value = 19.73722105032877
$ws.Cells.Item(25, 9).Value = 33.86301646839505
$ws.Cells.Item(25, 12).Value = 9.605077071469188
$ws.Cells.Item(25, 14).Value = 19.84905939529497

Write-Output "Updated loading_percent values for rows 2-25"
